# Refresh cryptos list values (prices / 1h volume %) for the daily
# GitHub Actions update, and swap the Kaspa/ARBITRUM rows (41/42).
# For Price (column D) cells whose new text looks like a plain number
# (e.g. "0.519", "0.0490"), the cell is pre-formatted as Text ("@") so
# COM stores the literal string instead of silently coercing it to a
# float (which would drop meaningful trailing zeros). The format is
# reset back to Normal afterwards so no stray number-format stays on
# the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.887.86"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "1.627.14"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.519"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.61"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.59%  "
$ws.Range("E9").Value = "  +3.40%  "
$ws.Range("E10").Value = "  +2.51%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "1.859.95"
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").Value = "1.627.14"
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("E14").Value = "  +6.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.64%  "
$ws.Range("D16").Value = "29.935.08"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("E17").Value = "  +19.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  +3.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.88%  "
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("E27").Value = "  +2.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.12%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0490"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.99%  "
$ws.Range("E31").Value = "  +6.56%  "
$ws.Range("E32").Value = "  +4.07%  "
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("D34").Value = "1.427.91"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("E35").Value = "  +6.52%  "
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("E37").Value = "  +2.17%  "
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("E39").Value = "  +3.30%  "
$ws.Range("E40").Value = "  +2.99%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.831"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.03%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0500"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "54.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.29%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "69.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.66%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.59%  "
$ws.Range("D49").Value = "1.768.12"
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "89.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.07%  "
$ws.Range("E51").Value = "  +0.76%  "
